$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh - update changed cells to match latest scrape.
# D-column price cells are forced to text via NumberFormat "@" so that values such as
# "218.30" or "0.999" are not silently reinterpreted as numbers (which would drop
# trailing zeros / change the stored type); the style is then reset back to "Normal"
# so no stray number-format style is left behind on the cell.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '92.750.58'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +4.63%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '3.271.51'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.17%  '
$ws.Range("E4").Value = '  +0.01%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '218.30'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.29%  '
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '629.33'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.402'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  +2.55%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.709'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +0.88%  '
$ws.Range("E9").Value = '  -0.04%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '3.266.47'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.15%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.588'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +1.81%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.0000270'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +2.41%  '
$ws.Range("E13").Value = '  -2.98%  '
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '34.24'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -0.04%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '92.340.21'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +4.53%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '3.862.89'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -0.31%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '5.33'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -0.36%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '3.275.30'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.31%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '3.28'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +4.78%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0000213'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  +59.33%  '
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '13.93'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -1.56%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '447.22'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +2.51%  '
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("E24").Value = '  -2.97%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '5.33'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +2.55%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '12.08'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -2.14%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '3.443.30'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +0.25%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '77.48'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +0.53%  '
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '0.173'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  -6.15%  '
$ws.Range("E31").Value = '  +0.04%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '8.73'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -1.77%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '552.99'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -3.48%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '3.83'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  +28.70%  '
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = '7.10'
$cell.Style = "Normal"
$ws.Range("E35").Value = '  +0.02%  '
$ws.Range("E36").Value = '  -1.94%  '
$ws.Range("E37").Value = '  -8.79%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '22.60'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  -0.28%  '
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '22.47'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +3.19%  '
$ws.Range("E40").Value = '  -7.04%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.998'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -0.07%  '
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '0.392'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  -1.57%  '
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '1.97'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -2.51%  '
$ws.Range("E44").Value = '  -0.05%  '
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '149.72'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -2.97%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '45.49'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  +1.37%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '178.64'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -0.92%  '
$ws.Range("E48").Value = '  +1.86%  '
$ws.Range("E49").Value = '  -1.21%  '
$ws.Range("E50").Value = '  +2.04%  '
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '4.21'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.70%  '
